$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the same "last changed" date serial for every
# data row (rows 2-344). Update it from 45206 to 45208 in one shot.
$ws.Range("C2:C344").Value = 45208
